# Apply the crypto price/volume refresh described by the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.184.65"
$ws.Range("E2").Value = "  -0.58%  "

$ws.Range("D3").Value = "1.862.77"
$ws.Range("E3").Value = "  -0.35%  "

$ws.Range("D4").Value = "'1.001"
$ws.Range("E4").Value = "  +0.03%  "

$ws.Range("D5").Value = "'0.7192"
$ws.Range("E5").Value = "  +1.16%  "

$ws.Range("D6").Value = "'240.61"
$ws.Range("E6").Value = "  +0.96%  "

$ws.Range("E7").Value = "  +0.14%  "

$ws.Range("D8").Value = "'0.07730"
$ws.Range("E8").Value = "  -0.84%  "

$ws.Range("D9").Value = "'0.3071"
$ws.Range("E9").Value = "  +0.25%  "

$ws.Range("D10").Value = "'24.95"
$ws.Range("E10").Value = "  -0.54%  "

$ws.Range("D11").Value = "'0.08255"
$ws.Range("E11").Value = "  +0.91%  "

$ws.Range("D12").Value = "1.886.96"
$ws.Range("E12").Value = "  +0.66%  "

# Row 13: Polygon/Polkadot swapped position in the ranking
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").Value = "'5.205"
$ws.Range("E13").Value = "  -0.75%  "

# Row 14: Polygon/Polkadot swapped position in the ranking
$ws.Range("B14").Value = "Polygon"
$ws.Range("C14").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D14").Value = "'0.7147"
$ws.Range("E14").Value = "  -0.67%  "

$ws.Range("D15").Value = "'90.10"
$ws.Range("E15").Value = "  +1.01%  "

$ws.Range("D16").Value = "29.210.48"
$ws.Range("E16").Value = "  -0.47%  "

$ws.Range("D17").Value = "'5.817"
$ws.Range("E17").Value = "  -0.12%  "

$ws.Range("D18").Value = "'242.74"
$ws.Range("E18").Value = "  +0.47%  "

$ws.Range("D19").Value = "'0.000007777"
$ws.Range("E19").Value = "  -0.47%  "

$ws.Range("D20").Value = "2.120.77"
$ws.Range("E20").Value = "  +0.59%  "

$ws.Range("D21").Value = "'13.09"
$ws.Range("E21").Value = "  -1.26%  "

$ws.Range("D22").Value = "'1.001"
$ws.Range("E22").Value = "  +0.19%  "

$ws.Range("D23").Value = "'7.941"
$ws.Range("E23").Value = "  +4.27%  "

$ws.Range("D24").Value = "'1.003"
$ws.Range("E24").Value = "  +0.12%  "

$ws.Range("D25").Value = "'0.1588"
$ws.Range("E25").Value = "  +9.40%  "

$ws.Range("D26").Value = "'162.03"
$ws.Range("E26").Value = "  -0.18%  "

$ws.Range("D27").Value = "'8.894"
$ws.Range("E27").Value = "  -0.70%  "

$ws.Range("D28").Value = "'18.15"
$ws.Range("E28").Value = "  -0.14%  "

$ws.Range("E29").Value = "  +0.71%  "

$ws.Range("D30").Value = "'1.297"
$ws.Range("E30").Value = "  -4.68%  "

$ws.Range("D31").Value = "'4.349"
$ws.Range("E31").Value = "  +0.47%  "

$ws.Range("D32").Value = "'4.081"
$ws.Range("E32").Value = "  +0.59%  "

$ws.Range("D33").Value = "'0.05176"
$ws.Range("E33").Value = "  -0.81%  "

$ws.Range("D34").Value = "'1.911"
$ws.Range("E34").Value = "  -0.46%  "

$ws.Range("D35").Value = "'1.173"
$ws.Range("E35").Value = "  -1.32%  "

$ws.Range("D36").Value = "'0.7276"
$ws.Range("E36").Value = "  +1.43%  "

$ws.Range("D37").Value = "'2.680"
$ws.Range("E37").Value = "  +0.18%  "

$ws.Range("D38").Value = "'0.01844"
$ws.Range("E38").Value = "  -0.38%  "

$ws.Range("D39").Value = "'2.700"
$ws.Range("E39").Value = "  -0.34%  "

$ws.Range("D40").Value = "1.157.08"
$ws.Range("E40").Value = "  -1.95%  "

$ws.Range("D41").Value = "'0.9000"
$ws.Range("E41").Value = "  -1.73%  "

$ws.Range("D42").Value = "'6.079"
$ws.Range("E42").Value = "  +1.19%  "

$ws.Range("D43").Value = "'72.02"
$ws.Range("E43").Value = "  +0.83%  "

$ws.Range("E44").Value = "  +0.17%  "

$ws.Range("D45").Value = "'101.59"
$ws.Range("E45").Value = "  -0.57%  "

$ws.Range("D46").Value = "2.015.36"
$ws.Range("E46").Value = "  +0.43%  "

$ws.Range("E47").Value = "  -1.41%  "

$ws.Range("D48").Value = "'1.761"
$ws.Range("E48").Value = "  -0.17%  "

$ws.Range("D49").Value = "'9.248"
$ws.Range("E49").Value = "  -0.04%  "

$ws.Range("D50").Value = "'2.859"
$ws.Range("E50").Value = "  +2.34%  "

$ws.Range("D51").Value = "'0.9982"
$ws.Range("E51").Value = "  -0.33%  "
